$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so numeric-looking
# strings (e.g. "2.00", "22.07") are not coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.852.31"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "2.748.44"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "577.42"
$ws.Range("E5").Value = "  -2.68%  "

$ws.Range("D6").Value = "158.95"
$ws.Range("E6").Value = "  +4.23%  "

$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").Value = "0.391"
$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("D11").Value = "5.69"
$ws.Range("E11").Value = "  -15.69%  "

$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "3.235.72"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "26.94"
$ws.Range("E14").Value = "  +0.82%  "

$ws.Range("D15").Value = "63.765.61"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "0.0000155"
$ws.Range("E16").Value = "  +1.22%  "

$ws.Range("D17").Value = "2.749.93"
$ws.Range("E17").Value = "  -0.70%  "

$ws.Range("D18").Value = "12.23"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "4.94"
$ws.Range("E19").Value = "  +0.77%  "

$ws.Range("D20").Value = "361.24"
$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("E21").Value = "  -2.06%  "

$ws.Range("D22").Value = "0.566"
$ws.Range("E22").Value = "  +5.31%  "

$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("D24").Value = "66.31"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("E25").Value = "  +2.76%  "

$ws.Range("D26").Value = "8.66"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").Value = "0.0₃0934"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  -1.69%  "

$ws.Range("D30").Value = "7.12"
$ws.Range("E30").Value = "  -0.27%  "

$ws.Range("E31").Value = "  +5.51%  "

$ws.Range("D32").Value = "168.72"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("E33").Value = "  +0.23%  "

$ws.Range("D34").Value = "20.54"
$ws.Range("E34").Value = "  -0.78%  "

$ws.Range("D35").Value = "5.01"
$ws.Range("E35").Value = "  +4.43%  "

$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  +2.73%  "

$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("E39").Value = "  -0.85%  "

$ws.Range("E40").Value = "  +8.97%  "

$ws.Range("D41").Value = "331.86"
$ws.Range("E41").Value = "  -4.93%  "

$ws.Range("D42").Value = "39.53"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").Value = "22.07"

$ws.Range("D44").Value = "0.0598"
$ws.Range("E44").Value = "  +0.71%  "

$ws.Range("D45").Value = "21.93"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.640"
$ws.Range("E46").Value = "  -1.64%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0258"
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").Value = "136.81"
$ws.Range("E48").Value = "  -5.03%  "

$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  +1.01%  "
